$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($r, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $values[$i]
    }
}

# Update the "last updated" timestamp banner
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 14 de Junio de 2020 a las 00:59"

# Refreshed case counts (values only change; row order stays the same except
# for the three row-pairs below that re-sort relative to their neighbour
# because the refreshed totals changed their rank).
Set-Row 4   @("Estados Unidos", 2140992, 24070, 849355, 1174135, 0, 677, 117502)
Set-Row 7   @("India", 321626, 12023, 162326, 150101, 0, 309, 9199)
Set-Row 20  @("Canada", 98410, 467, 59353, 30950, 0, 58, 8107)

# Colombia overtakes Paises Bajos
Set-Row 28  @("Colombia", 48746, 1888, 19426, 27728, 0, 47, 1592)
Set-Row 29  @("Paises Bajos", 48640, 179, 0, 0, 0, 4, 6057)

# Argentina overtakes Polonia
Set-Row 39  @("Argentina", 30295, 1531, 9083, 20397, 0, 30, 815)
Set-Row 40  @("Polonia", 29017, 440, 14104, 13676, 0, 15, 1237)

Set-Row 51  @("Japon", 17382, 50, 15580, 878, 0, 2, 924)
Set-Row 55  @("Nigeria", 15682, 501, 5101, 10174, 0, 8, 407)
Set-Row 68  @("Noruega", 8628, 8, 8138, 248, 0, 0, 242)
Set-Row 127 @("Niger", 980, 2, 881, 33, 0, 1, 66)
Set-Row 133 @("Uruguay", 847, 0, 784, 40, 0, 0, 23)
Set-Row 163 @("Surinam", 196, 9, 9, 184, 0, 0, 3)

# Montserrat overtakes Seychelles
Set-Row 210 @("Montserrat", 11, 0, 10, 0, 0, 0, 1)
Set-Row 211 @("Seychelles", 11, 0, 11, 0, 0, 0, 0)
